$d = $word.ActiveDocument

# The table's sole section was portrait (pgSz w:h=16848 w:w=11952); the
# commit flips it to landscape (pgSz w:h=11952 w:w=16848, orient=landscape)
# so the wide results table fits the page. wdOrientLandscape = 1.
$d.PageSetup.Orientation = 1
